# Auto-generated edit script: refresh market-price-derived columns (H,I,J,K,L,M,N)
# on the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1307.8
$ws.Range("I5").Value = 1253.8
$ws.Range("K5").Value = 1253.8
$ws.Range("M5").Value = -1138.8
$ws.Range("H52").Value = 3345
$ws.Range("I52").Value = 2283
$ws.Range("K52").Value = 6849
$ws.Range("M52").Value = -6689
$ws.Range("H62").Value = 4170.294
$ws.Range("I62").Value = 4381.4546
$ws.Range("K62").Value = 4381.4546
$ws.Range("M62").Value = -3757.4546
$ws.Range("H64").Value = 6899.4287
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 4170.294
$ws.Range("I65").Value = 4381.4546
$ws.Range("K65").Value = 21907.273
$ws.Range("M65").Value = -18787.273
$ws.Range("H67").Value = 6899.4287
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H86").Value = 29115
$ws.Range("I86").Value = 761.2
$ws.Range("K86").Value = 761.2
$ws.Range("M86").Value = 361.8
$ws.Range("H89").Value = 29115
$ws.Range("I89").Value = 761.2
$ws.Range("K89").Value = 3806
$ws.Range("M89").Value = 1810
$ws.Range("H106").Value = 12209.314
$ws.Range("I106").Value = 13517.419
$ws.Range("K106").Value = 13517.419
$ws.Range("M106").Value = -12886.419
$ws.Range("H121").Value = 1817.125
$ws.Range("J121").Value = 1817.125
$ws.Range("L121").Value = 5451.375
$ws.Range("N121").Value = -8945.375
$ws.Range("H126").Value = 99975.5
$ws.Range("J126").Value = 99975.5
$ws.Range("L126").Value = 99975.5
$ws.Range("N126").Value = -109855.5
$ws.Range("H138").Value = 4724.9565
$ws.Range("I138").Value = 2071.6875
$ws.Range("K138").Value = 6215.0625
$ws.Range("M138").Value = -1075.0625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2424
$ws.Range("I88").Value = 1999.5
$ws.Range("K88").Value = 1999.5
$ws.Range("M88").Value = -1593.5
$ws.Range("H91").Value = 2424
$ws.Range("I91").Value = 1999.5
$ws.Range("K91").Value = 1999.5
$ws.Range("M91").Value = -595.5
$ws.Range("H114").Value = 99998.5
$ws.Range("J114").Value = 99998.5
$ws.Range("L114").Value = 99998.5
$ws.Range("N114").Value = -108676.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 13588.25
$ws.Range("I96").Value = 13588.25
$ws.Range("K96").Value = 13588.25
$ws.Range("M96").Value = -10842.25
$ws.Range("H134").Value = 2513.2856
$ws.Range("I134").Value = 2098.8333
$ws.Range("K134").Value = 6296.499899999999
$ws.Range("M134").Value = -3761.499899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3422.923
$ws.Range("I22").Value = 3700
$ws.Range("J22").Value = 2799.5
$ws.Range("K22").Value = 3700
$ws.Range("L22").Value = 2799.5
$ws.Range("M22").Value = -3350
$ws.Range("N22").Value = -3499.5
$ws.Range("H62").Value = 3901.1667
$ws.Range("I62").Value = 3673
$ws.Range("J62").Value = 4699.75
$ws.Range("K62").Value = 3673
$ws.Range("L62").Value = 4699.75
$ws.Range("M62").Value = -3049
$ws.Range("N62").Value = -5947.75
$ws.Range("H65").Value = 3901.1667
$ws.Range("I65").Value = 3673
$ws.Range("J65").Value = 4699.75
$ws.Range("K65").Value = 18365
$ws.Range("L65").Value = 23498.75
$ws.Range("M65").Value = -15245
$ws.Range("N65").Value = -29738.75
$ws.Range("H134").Value = 1987.3948
$ws.Range("I134").Value = 1895.1936
$ws.Range("K134").Value = 5685.5808
$ws.Range("M134").Value = -3150.5808
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 460.27274
$ws.Range("J97").Value = 525.26666
$ws.Range("L97").Value = 1575.79998
$ws.Range("N97").Value = -2567.79998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 10373.25
$ws.Range("J59").Value = 5749.5
$ws.Range("L59").Value = 5749.5
$ws.Range("N59").Value = -6915.5
$ws.Range("H70").Value = 29697.62
$ws.Range("I70").Value = 41622.645
$ws.Range("K70").Value = 41622.645
$ws.Range("M70").Value = -41352.645
$ws.Range("H73").Value = 29697.62
$ws.Range("I73").Value = 41622.645
$ws.Range("K73").Value = 41622.645
$ws.Range("M73").Value = -40686.645
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7784.6665
$ws.Range("I40").Value = 7741.6
$ws.Range("K40").Value = 7741.6
$ws.Range("M40").Value = -7605.6
$ws.Range("H46").Value = 992.7143
$ws.Range("J46").Value = 993.2
$ws.Range("L46").Value = 993.2
$ws.Range("N46").Value = -1369.2
$ws.Range("H68").Value = 3437.2778
$ws.Range("I68").Value = 1812.6923
$ws.Range("K68").Value = 1812.6923
$ws.Range("M68").Value = -1063.6923
$ws.Range("H71").Value = 3437.2778
$ws.Range("I71").Value = 1812.6923
$ws.Range("K71").Value = 9063.4615
$ws.Range("M71").Value = -5319.461499999999
$ws.Range("H122").Value = 9163.333
$ws.Range("I122").Value = 9163.333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 27489.999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -25039.999
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 5504.6875
$ws.Range("I136").Value = 5007.222
$ws.Range("K136").Value = 15021.666
$ws.Range("M136").Value = -12471.666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H62").Value = 550000
$ws.Range("I62").Value = 550000
$ws.Range("K62").Value = 550000
$ws.Range("M62").Value = -549376
$ws.Range("H65").Value = 550000
$ws.Range("I65").Value = 550000
$ws.Range("K65").Value = 2750000
$ws.Range("M65").Value = -2746880
$ws.Range("H81").Value = 4304.727
$ws.Range("I81").Value = 2396.3333
$ws.Range("K81").Value = 4792.6666
$ws.Range("M81").Value = -3731.6666
$ws.Range("H84").Value = 4304.727
$ws.Range("I84").Value = 2396.3333
$ws.Range("K84").Value = 23963.333
$ws.Range("M84").Value = -18659.333
$ws.Range("H126").Value = 9527626
$ws.Range("I126").Value = 15876740
$ws.Range("K126").Value = 47630220
$ws.Range("M126").Value = -47627750

Write-Host "Applied Zodiark_Profits market data refresh."
